$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(58, 1).Value = '1b650324aafac0e6e6f0c473eccff258'
$ws.Cells.Item(58, 2).Value = 'Homo sapiens'
$ws.Cells.Item(58, 3).Value = 'Human'
$ws.Cells.Item(58, 4).Value = 'Human'
$ws.Cells.Item(59, 1).Value = '407f08a29007a8a153222d82ef47d408'
$ws.Cells.Item(59, 2).Value = 'Menidia menidia'
$ws.Cells.Item(59, 3).Value = 'Atlantic silverside'
$ws.Cells.Item(59, 4).Value = 'Teleost Fish'
$ws.Cells.Item(76, 1).Value = '09351b480c58a99c4150d54ebbc97c6c'
$ws.Cells.Item(76, 2).Value = 'Menidia menidia'
$ws.Cells.Item(76, 3).Value = 'Atlantic silverside'
$ws.Cells.Item(76, 4).Value = 'Teleost Fish'
$ws.Cells.Item(77, 1).Value = 'd34820a8c9954e292ea9dbc76f4275b4'
$ws.Cells.Item(77, 2).Value = 'Homo sapiens'
$ws.Cells.Item(77, 3).Value = 'Human'
$ws.Cells.Item(77, 4).Value = 'Human'
$ws.Cells.Item(93, 1).Value = '680475954df3011ebba1033f1b2f2a86'
$ws.Cells.Item(93, 2).Value = 'Prionotus carolinus'
$ws.Cells.Item(93, 3).Value = 'Northern sea robin'
$ws.Cells.Item(93, 4).Value = 'Teleost Fish'
$ws.Cells.Item(94, 1).Value = '4c5905c5ab539613d9c3069d0ae54188'
$ws.Cells.Item(94, 2).Value = 'Bos taurus'
$ws.Cells.Item(94, 3).Value = 'Cow'
$ws.Cells.Item(94, 4).Value = 'Livestock'
$ws.Cells.Item(95, 1).Value = 'db8615250f29272019fe417d96bf08f3'
$ws.Cells.Item(95, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(95, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(95, 4).Value = 'Teleost Fish'
$ws.Cells.Item(105, 1).Value = 'f5e0ea6fe3e45da9605b758c440ae692'
$ws.Cells.Item(105, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(105, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(105, 4).Value = 'Teleost Fish'
$ws.Cells.Item(106, 1).Value = '6d1668646cf923fa90217b0797de7a7d'
$ws.Cells.Item(106, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(106, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(106, 4).Value = 'Teleost Fish'
$ws.Cells.Item(113, 1).Value = '5e733a21f67e541f28ed4bf4fe025044'
$ws.Cells.Item(113, 2).Value = 'Paralichthys dentatus'
$ws.Cells.Item(113, 3).Value = 'Summer flounder'
$ws.Cells.Item(113, 4).Value = 'Teleost Fish'
$ws.Cells.Item(114, 1).Value = 'f524c4b860dec1e6b994c28dd8e4b75e'
$ws.Cells.Item(114, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(114, 3).Value = 'Northern sand lance'
$ws.Cells.Item(114, 4).Value = 'Teleost Fish'
$ws.Cells.Item(144, 1).Value = '9db3dc01519672b43908456a37b27b4d'
$ws.Cells.Item(144, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(144, 3).Value = 'Mummichog'
$ws.Cells.Item(144, 4).Value = 'Teleost Fish'
$ws.Cells.Item(145, 1).Value = '0df37a1b74088f0e7410a1b78cada881'
$ws.Cells.Item(145, 2).Value = 'Engraulis eurystole'
$ws.Cells.Item(145, 3).Value = 'Silver anchovy'
$ws.Cells.Item(145, 4).Value = 'Teleost Fish'
$ws.Cells.Item(150, 1).Value = '0f011be680aec3ee4b12b1b139902251'
$ws.Cells.Item(150, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(150, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(150, 4).Value = 'Teleost Fish'
$ws.Cells.Item(151, 1).Value = '94a944154183c458facbab20fe39ffa9'
$ws.Cells.Item(151, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(151, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(151, 4).Value = 'Teleost Fish'
$ws.Cells.Item(172, 1).Value = '191ed810bb884ed43fa1919f6da3d82a'
$ws.Cells.Item(172, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(172, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(172, 4).Value = 'Teleost Fish'
$ws.Cells.Item(173, 1).Value = '148aa3594130e12c353383f68bfa0b6a'
$ws.Cells.Item(173, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(173, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(173, 4).Value = 'Teleost Fish'
$ws.Cells.Item(174, 1).Value = 'c73cefb2b4ac8de08ae0c68341cbb28f'
$ws.Cells.Item(174, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(174, 3).Value = 'Northern sand lance'
$ws.Cells.Item(174, 4).Value = 'Teleost Fish'
$ws.Cells.Item(175, 1).Value = 'ed00c11476e9a07d3441cb0a1073d3ab'
$ws.Cells.Item(175, 2).Value = 'Etropus microstomus'
$ws.Cells.Item(175, 3).Value = 'Smallmouth flounder'
$ws.Cells.Item(175, 4).Value = 'Teleost Fish'
$ws.Cells.Item(181, 1).Value = 'ff405ebc8992c59ba51a99e33a12fe74'
$ws.Cells.Item(181, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(181, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(181, 4).Value = 'Teleost Fish'
$ws.Cells.Item(182, 1).Value = '558d8758ae62abe36b1507ce2094ef7c'
$ws.Cells.Item(182, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(182, 3).Value = 'Northern sand lance'
$ws.Cells.Item(182, 4).Value = 'Teleost Fish'
$ws.Cells.Item(186, 1).Value = 'ee3c408644b66e62dde706ff463f359a'
$ws.Cells.Item(186, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(186, 3).Value = 'Northern sand lance'
$ws.Cells.Item(186, 4).Value = 'Teleost Fish'
$ws.Cells.Item(187, 1).Value = '86c340102750abe5f2a75f3d5501b55d'
$ws.Cells.Item(187, 2).Value = 'Menidia beryllina'
$ws.Cells.Item(187, 3).Value = 'Inland silverside'
$ws.Cells.Item(187, 4).Value = 'Teleost Fish'
$ws.Cells.Item(190, 1).Value = '0ad9142dc74ab0ef2021cfff48d4194d'
$ws.Cells.Item(190, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(190, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(190, 4).Value = 'Teleost Fish'
$ws.Cells.Item(191, 1).Value = '8830d0cf4452e1cd0f9a6552b48b2b40'
$ws.Cells.Item(191, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(191, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(191, 4).Value = 'Teleost Fish'
$ws.Cells.Item(192, 1).Value = '731abf4fa491ab03dd796729de5ab3eb'
$ws.Cells.Item(192, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(192, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(192, 4).Value = 'Teleost Fish'
$ws.Cells.Item(195, 1).Value = 'f2e15a0b398b704a888c965d3b49035b'
$ws.Cells.Item(195, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(195, 3).Value = 'Mummichog'
$ws.Cells.Item(195, 4).Value = 'Teleost Fish'
$ws.Cells.Item(205, 1).Value = 'e468b57f39f048ada7562924022dc516'
$ws.Cells.Item(205, 2).Value = 'Homo sapiens'
$ws.Cells.Item(205, 3).Value = 'Human'
$ws.Cells.Item(205, 4).Value = 'Human'
$ws.Cells.Item(206, 1).Value = '1533469db84e906a7d07208d202f0b61'
$ws.Cells.Item(206, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(206, 3).Value = 'Mummichog'
$ws.Cells.Item(206, 4).Value = 'Teleost Fish'
$ws.Cells.Item(216, 1).Value = '5b2278535af7a77c15966bc43d0188bd'
$ws.Cells.Item(216, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(216, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(216, 4).Value = 'Teleost Fish'
$ws.Cells.Item(217, 1).Value = '6a83eb23e34e01773abb7d038e38c583'
$ws.Cells.Item(217, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(217, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(217, 4).Value = 'Teleost Fish'
$ws.Cells.Item(218, 1).Value = 'bdb87097756f45aa57e56f1d9f456f26'
$ws.Cells.Item(218, 2).Value = 'Larus sp'
$ws.Cells.Item(218, 3).Value = 'Great black backed gull and other gulls'
$ws.Cells.Item(218, 4).Value = 'Bird'
$ws.Cells.Item(220, 1).Value = 'cb17be39fabe38eb2368ba0635321393'
$ws.Cells.Item(220, 2).Value = 'Tautoga onitis'
$ws.Cells.Item(220, 3).Value = 'Tautog'
$ws.Cells.Item(220, 4).Value = 'Teleost Fish'
$ws.Cells.Item(221, 1).Value = 'c0b18824ab60460cd31eed51f737f882'
$ws.Cells.Item(221, 2).Value = 'Homo sapiens'
$ws.Cells.Item(221, 3).Value = 'Human'
$ws.Cells.Item(221, 4).Value = 'Human'
$ws.Cells.Item(225, 1).Value = '879319f127f42872ba2daeb54fc4135a'
$ws.Cells.Item(225, 2).Value = 'Gasterosteus aculeatus'
$ws.Cells.Item(225, 3).Value = 'Threespined stickleback'
$ws.Cells.Item(225, 4).Value = 'Teleost Fish'
$ws.Cells.Item(226, 1).Value = 'f937641d91db232cb7180be9e04fb9e0'
$ws.Cells.Item(226, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(226, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(226, 4).Value = 'Teleost Fish'
$ws.Cells.Item(227, 1).Value = '977b02be79d865979e54848db649eaf0'
$ws.Cells.Item(227, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(227, 3).Value = 'Northern sand lance'
$ws.Cells.Item(227, 4).Value = 'Teleost Fish'
$ws.Cells.Item(228, 1).Value = 'df263dae379496c7e522db8a7dbc01c9'
$ws.Cells.Item(228, 2).Value = 'Scomber scombrus'
$ws.Cells.Item(228, 3).Value = 'Atlantic mackerel'
$ws.Cells.Item(228, 4).Value = 'Teleost Fish'
$ws.Cells.Item(229, 1).Value = '4450a6fa10b56881617cff33c5585aa8'
$ws.Cells.Item(229, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(229, 3).Value = 'Mummichog'
$ws.Cells.Item(229, 4).Value = 'Teleost Fish'
$ws.Cells.Item(230, 1).Value = 'f5ca5d430f1b145903b92fc335a4bafd'
$ws.Cells.Item(230, 2).Value = 'Lucania parva'
$ws.Cells.Item(230, 3).Value = 'Rainwater killifish'
$ws.Cells.Item(230, 4).Value = 'Teleost Fish'
$ws.Cells.Item(231, 1).Value = '9c8a7b893d0fdaf8c1c89606cfce1c08'
$ws.Cells.Item(231, 2).Value = 'Enchelyopus cimbrius'
$ws.Cells.Item(231, 3).Value = 'Fourbeard rockling'
$ws.Cells.Item(231, 4).Value = 'Teleost Fish'
$ws.Cells.Item(232, 1).Value = 'c73c11c0d8b73d825099e12aaaefb637'
$ws.Cells.Item(232, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(232, 3).Value = 'Mummichog'
$ws.Cells.Item(232, 4).Value = 'Teleost Fish'
$ws.Cells.Item(233, 1).Value = 'c972de9c10572043855aaca4a4da68f4'
$ws.Cells.Item(233, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(233, 3).Value = 'Northern sand lance'
$ws.Cells.Item(233, 4).Value = 'Teleost Fish'
$ws.Cells.Item(234, 1).Value = '243cc9259a8d104346a5dd517ca99499'
$ws.Cells.Item(234, 2).Value = 'Ammodytes americanus'
$ws.Cells.Item(234, 3).Value = 'American sand lance'
$ws.Cells.Item(234, 4).Value = 'Teleost Fish'
$ws.Cells.Item(244, 1).Value = 'f753730afbaa726c79bd991f32ea9778'
$ws.Cells.Item(244, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(244, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(244, 4).Value = 'Teleost Fish'
$ws.Cells.Item(245, 1).Value = '5432a6e652c21bb79c110c1179832080'
$ws.Cells.Item(245, 2).Value = 'Clupeidae sp'
$ws.Cells.Item(245, 3).Value = 'Atlantic menhaden or River herrings'
$ws.Cells.Item(245, 4).Value = 'Teleost Fish'
$ws.Cells.Item(246, 1).Value = '279fde05b5aed4bfe15ab39776ff82ba'
$ws.Cells.Item(246, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(246, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(246, 4).Value = 'Teleost Fish'
$ws.Cells.Item(248, 1).Value = '4db280926cca07cc86b0e098513d9cc0'
$ws.Cells.Item(248, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(248, 3).Value = 'Northern sand lance'
$ws.Cells.Item(248, 4).Value = 'Teleost Fish'
$ws.Cells.Item(249, 1).Value = '88065f0fd14ae3b76fc1a87f8df6ef2d'
$ws.Cells.Item(249, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(249, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(249, 4).Value = 'Teleost Fish'
$ws.Cells.Item(251, 1).Value = 'fc6d040e1564a91e1c6d67e1e32b9022'
$ws.Cells.Item(251, 2).Value = 'Rattus norvegicus'
$ws.Cells.Item(251, 3).Value = 'Norway rat'
$ws.Cells.Item(251, 4).Value = 'Mammal'
$ws.Cells.Item(258, 1).Value = '14bd3bb11b9a6c641ad60556bf6141d0'
$ws.Cells.Item(258, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(258, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(258, 4).Value = 'Teleost Fish'
$ws.Cells.Item(259, 1).Value = '8103469b2716037f1cc4ce8959ae0081'
$ws.Cells.Item(259, 2).Value = 'Menidia menidia'
$ws.Cells.Item(259, 3).Value = 'Atlantic silverside'
$ws.Cells.Item(259, 4).Value = 'Teleost Fish'
$ws.Cells.Item(262, 1).Value = 'f4d5447013c09b659b99d47459de2042'
$ws.Cells.Item(262, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(262, 3).Value = 'Mummichog'
$ws.Cells.Item(262, 4).Value = 'Teleost Fish'
$ws.Cells.Item(263, 1).Value = '29d8e064f48ae7211c9fba32872b36f9'
$ws.Cells.Item(263, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(263, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(263, 4).Value = 'Teleost Fish'
$ws.Cells.Item(264, 1).Value = '53cfac0a209f1dbdaf758a75c84df7d6'
$ws.Cells.Item(264, 2).Value = 'Homo sapiens'
$ws.Cells.Item(264, 3).Value = 'Human'
$ws.Cells.Item(264, 4).Value = 'Human'
$ws.Cells.Item(265, 1).Value = '9dc70d19a67c006232234c9bcbbab33f'
$ws.Cells.Item(265, 2).Value = 'Clangula hyemalis or other Anatidae sp'
$ws.Cells.Item(265, 3).Value = 'Long tailed duck or other ducks'
$ws.Cells.Item(265, 4).Value = 'Bird'
$ws.Cells.Item(267, 1).Value = 'daa05108bed6292fbd2eedef6214fdff'
$ws.Cells.Item(267, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(267, 3).Value = 'Mummichog'
$ws.Cells.Item(267, 4).Value = 'Teleost Fish'
$ws.Cells.Item(269, 1).Value = '4c451c35f94e737edb8816211106c35d'
$ws.Cells.Item(269, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(269, 3).Value = 'Mummichog'
$ws.Cells.Item(269, 4).Value = 'Teleost Fish'
$ws.Cells.Item(270, 1).Value = 'b61cc19540f9627af5070110979ebf91'
$ws.Cells.Item(270, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(270, 3).Value = 'Northern sand lance'
$ws.Cells.Item(270, 4).Value = 'Teleost Fish'
$ws.Cells.Item(271, 1).Value = '26a933bf07de2306f33cc95ec94e4b2f'
$ws.Cells.Item(271, 2).Value = 'Pomoxis nigromaculatus'
$ws.Cells.Item(271, 3).Value = 'Black crappie'
$ws.Cells.Item(271, 4).Value = 'Teleost Fish'
$ws.Cells.Item(273, 1).Value = '93b36a6e82074114f7c4d90b6172dba2'
$ws.Cells.Item(273, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(273, 3).Value = 'Mummichog'
$ws.Cells.Item(273, 4).Value = 'Teleost Fish'
$ws.Cells.Item(274, 1).Value = '1be2860881962b4dd3a0a7c6db14ca80'
$ws.Cells.Item(274, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(274, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(274, 4).Value = 'Teleost Fish'
$ws.Cells.Item(278, 1).Value = 'b7b35bf53a25eef31602b3f785c925e9'
$ws.Cells.Item(278, 2).Value = 'Anguilla rostrata'
$ws.Cells.Item(278, 3).Value = 'American eel'
$ws.Cells.Item(278, 4).Value = 'Teleost Fish'
$ws.Cells.Item(280, 1).Value = '5cc8cf140b434a5aba3a1dad41339918'
$ws.Cells.Item(280, 2).Value = 'Anguilla rostrata'
$ws.Cells.Item(280, 3).Value = 'American eel'
$ws.Cells.Item(280, 4).Value = 'Teleost Fish'
$ws.Cells.Item(283, 1).Value = '29df87a23b45339e281dc7b390d16860'
$ws.Cells.Item(283, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(283, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(283, 4).Value = 'Teleost Fish'
$ws.Cells.Item(284, 1).Value = '8dde047966dfd43b699a5ca7122e55d7'
$ws.Cells.Item(284, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(284, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(284, 4).Value = 'Teleost Fish'
$ws.Cells.Item(285, 1).Value = '118de0da9053ad27ad0e3c1e136454d9'
$ws.Cells.Item(285, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(285, 3).Value = 'Northern sand lance'
$ws.Cells.Item(285, 4).Value = 'Teleost Fish'
$ws.Cells.Item(287, 1).Value = 'ed57094384d61f9a2dbc0c0e4ff6fb8a'
$ws.Cells.Item(287, 2).Value = 'Lepomis gibbosus'
$ws.Cells.Item(287, 3).Value = 'Pumpkinseed'
$ws.Cells.Item(287, 4).Value = 'Teleost Fish'
$ws.Cells.Item(288, 1).Value = 'd8de668ed4c19b2b2ea845f3db18ae79'
$ws.Cells.Item(288, 2).Value = 'Cottidae sp'
$ws.Cells.Item(288, 3).Value = 'Sculpins'
$ws.Cells.Item(288, 4).Value = 'Teleost Fish'
$ws.Cells.Item(289, 1).Value = '7bb99841c9c9b86e238cbcddc0e16567'
$ws.Cells.Item(289, 2).Value = 'Micropterus salmoides'
$ws.Cells.Item(289, 3).Value = 'Largemouth bass'
$ws.Cells.Item(289, 4).Value = 'Teleost Fish'
$ws.Cells.Item(292, 1).Value = '1dd0f1ca2adf649d8cba813ea6e43de2'
$ws.Cells.Item(292, 2).Value = 'Paralichthys dentatus'
$ws.Cells.Item(292, 3).Value = 'Summer flounder'
$ws.Cells.Item(292, 4).Value = 'Teleost Fish'
$ws.Cells.Item(293, 1).Value = '50bddde558bebcd7fa8dbf6542ea44fe'
$ws.Cells.Item(293, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(293, 3).Value = 'Mummichog'
$ws.Cells.Item(293, 4).Value = 'Teleost Fish'
$ws.Cells.Item(295, 1).Value = '92693323f831e69117617606814ae81f'
$ws.Cells.Item(295, 2).Value = 'Coryphaena hippurus'
$ws.Cells.Item(295, 3).Value = 'Mahi mahi'
$ws.Cells.Item(295, 4).Value = 'Teleost Fish'
$ws.Cells.Item(297, 1).Value = '0d6e610cd1019f50d693803e46db364f'
$ws.Cells.Item(297, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(297, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(297, 4).Value = 'Teleost Fish'
$ws.Cells.Item(298, 1).Value = 'bbd6b723329db44753870a56d15bdbd6'
$ws.Cells.Item(298, 2).Value = 'Sebastes fasciatus'
$ws.Cells.Item(298, 3).Value = 'Acadian redfish'
$ws.Cells.Item(298, 4).Value = 'Teleost Fish'
$ws.Cells.Item(303, 1).Value = '9f4285ab8775db6b862ee4fb416f0f5d'
$ws.Cells.Item(303, 2).Value = 'Melospiza melodia or Spizella passerina'
$ws.Cells.Item(303, 3).Value = 'Song sparrow or Chipping sparrow'
$ws.Cells.Item(303, 4).Value = 'Bird'
$ws.Cells.Item(304, 1).Value = '0a6108b3c6bbca90164970efbea23261'
$ws.Cells.Item(304, 2).Value = 'Apeltes quadracus'
$ws.Cells.Item(304, 3).Value = 'Fourspine stickleback'
$ws.Cells.Item(304, 4).Value = 'Teleost Fish'
$ws.Cells.Item(305, 1).Value = '901fc1f68af659cc3f6678c6a7396845'
$ws.Cells.Item(305, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(305, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(305, 4).Value = 'Teleost Fish'
$ws.Cells.Item(307, 1).Value = '5d4b77f374dcda6b5f48e88cc2b9664b'
$ws.Cells.Item(307, 2).Value = 'Unassigned'
$ws.Cells.Item(307, 3).Value = 'Unassigned'
$ws.Cells.Item(307, 4).Value = 'Unassigned'
$ws.Cells.Item(308, 1).Value = '7950b1078efc076defba9c936b970ef7'
$ws.Cells.Item(308, 2).Value = 'Homo sapiens'
$ws.Cells.Item(308, 3).Value = 'Human'
$ws.Cells.Item(308, 4).Value = 'Human'
$ws.Cells.Item(309, 1).Value = '8ea2a9236bef33ba65acfc82e6947942'
$ws.Cells.Item(309, 2).Value = 'Sciurus carolinensis'
$ws.Cells.Item(309, 3).Value = 'Gray squirrel'
$ws.Cells.Item(309, 4).Value = 'Mammal'
$ws.Cells.Item(310, 1).Value = 'c1f17b3dc22ac71ee83288f654c93bb3'
$ws.Cells.Item(310, 2).Value = 'Lontra canadensis'
$ws.Cells.Item(310, 3).Value = 'River otter'
$ws.Cells.Item(310, 4).Value = 'Mammal'
$ws.Cells.Item(311, 1).Value = 'a4e8997c6347c55b72f81e0accce0c37'
$ws.Cells.Item(311, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(311, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(311, 4).Value = 'Teleost Fish'
$ws.Cells.Item(313, 1).Value = '7eef5797ad87b51600785f22606c70de'
$ws.Cells.Item(313, 2).Value = 'Brevoortia tyrannus'
$ws.Cells.Item(313, 3).Value = 'Atlantic menhaden'
$ws.Cells.Item(313, 4).Value = 'Teleost Fish'
$ws.Cells.Item(314, 1).Value = '033531a8711295f5cf38c1111629eb77'
$ws.Cells.Item(314, 2).Value = 'Anguilla rostrata'
$ws.Cells.Item(314, 3).Value = 'American eel'
$ws.Cells.Item(314, 4).Value = 'Teleost Fish'
$ws.Cells.Item(316, 1).Value = '16d55edf1062cb60bf8a36a1da3212b5'
$ws.Cells.Item(316, 2).Value = 'Esox americanus or niger'
$ws.Cells.Item(316, 3).Value = 'Grass or chain pickerel'
$ws.Cells.Item(316, 4).Value = 'Teleost Fish'
$ws.Cells.Item(317, 1).Value = '5b1dbdcc719bcfd9ea209ec7d9ecd075'
$ws.Cells.Item(317, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(317, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(317, 4).Value = 'Teleost Fish'
$ws.Cells.Item(319, 1).Value = 'd1de955bd9480b3f0c70f78caec5a443'
$ws.Cells.Item(319, 2).Value = 'Homo sapiens'
$ws.Cells.Item(319, 3).Value = 'Human'
$ws.Cells.Item(319, 4).Value = 'Human'
$ws.Cells.Item(320, 1).Value = 'bfab25a003878187c8038ee55fdb7a53'
$ws.Cells.Item(320, 2).Value = 'Gavia immer'
$ws.Cells.Item(320, 3).Value = 'Common loon'
$ws.Cells.Item(320, 4).Value = 'Bird'
$ws.Cells.Item(321, 1).Value = 'fb3bb0a4483dcfbc39e8b7ccf8196749'
$ws.Cells.Item(321, 2).Value = 'Homo sapiens'
$ws.Cells.Item(321, 3).Value = 'Human'
$ws.Cells.Item(321, 4).Value = 'Human'
$ws.Cells.Item(322, 1).Value = '9ed3306f1d1dfb81749820128e325abc'
$ws.Cells.Item(322, 2).Value = 'Clupeidae sp'
$ws.Cells.Item(322, 3).Value = 'Atlantic menhaden or River herrings'
$ws.Cells.Item(322, 4).Value = 'Teleost Fish'
$ws.Cells.Item(327, 1).Value = '9f2355fd161fec0177a83045e771a239'
$ws.Cells.Item(327, 2).Value = 'Decapterus punctatus'
$ws.Cells.Item(327, 3).Value = 'Round scad'
$ws.Cells.Item(327, 4).Value = 'Teleost Fish'
$ws.Cells.Item(328, 1).Value = '032747b5e01bbdc74f3e3e59d9c2275f'
$ws.Cells.Item(328, 2).Value = 'Unassigned'
$ws.Cells.Item(328, 3).Value = 'Unassigned'
$ws.Cells.Item(328, 4).Value = 'Unassigned'
$ws.Cells.Item(336, 1).Value = '856a99622e9c49fc86cdcf3ddbefcfd1'
$ws.Cells.Item(336, 2).Value = 'Gobiosoma ginsburgi'
$ws.Cells.Item(336, 3).Value = 'Seaboard goby'
$ws.Cells.Item(336, 4).Value = 'Teleost Fish'
$ws.Cells.Item(337, 1).Value = '38e4381a618398d035b19d8c47e8bfa3'
$ws.Cells.Item(337, 2).Value = 'Ammodytes dubius'
$ws.Cells.Item(337, 3).Value = 'Northern sand lance'
$ws.Cells.Item(337, 4).Value = 'Teleost Fish'
$ws.Cells.Item(338, 1).Value = 'a1f66fcba0bd12eed7a60901bd8c6010'
$ws.Cells.Item(338, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(338, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(338, 4).Value = 'Teleost Fish'
$ws.Cells.Item(340, 1).Value = '5cde257b7febb75f7c9848a21bfe18cd'
$ws.Cells.Item(340, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(340, 3).Value = 'Mummichog'
$ws.Cells.Item(340, 4).Value = 'Teleost Fish'
$ws.Cells.Item(341, 1).Value = 'd3b57d4cf93def7c41d4b9baced940d9'
$ws.Cells.Item(341, 2).Value = 'Canis lupus'
$ws.Cells.Item(341, 3).Value = 'Dog'
$ws.Cells.Item(341, 4).Value = 'Livestock'
$ws.Cells.Item(343, 1).Value = '60fb36f888cbf7e4639c1bb98f0adc57'
$ws.Cells.Item(343, 2).Value = 'Menidia menidia'
$ws.Cells.Item(343, 3).Value = 'Atlantic silverside'
$ws.Cells.Item(343, 4).Value = 'Teleost Fish'
$ws.Cells.Item(345, 1).Value = 'e7f90ff8c7b97da66ce6d940d857e8e4'
$ws.Cells.Item(345, 2).Value = 'Trachurus lathami'
$ws.Cells.Item(345, 3).Value = 'Rough scad'
$ws.Cells.Item(345, 4).Value = 'Teleost Fish'
$ws.Cells.Item(346, 1).Value = 'bcf13bc540e00c02358754d8a1b40a9c'
$ws.Cells.Item(346, 2).Value = 'Fundulus heteroclitus'
$ws.Cells.Item(346, 3).Value = 'Mummichog'
$ws.Cells.Item(346, 4).Value = 'Teleost Fish'
$ws.Cells.Item(347, 1).Value = '7fef2f8e6a8bee56528216dfc05f0d81'
$ws.Cells.Item(347, 2).Value = 'Homo sapiens'
$ws.Cells.Item(347, 3).Value = 'Human'
$ws.Cells.Item(347, 4).Value = 'Human'
$ws.Cells.Item(353, 1).Value = '1a9a786e8451eec71300762a5398f4d0'
$ws.Cells.Item(353, 2).Value = 'Anguilla rostrata'
$ws.Cells.Item(353, 3).Value = 'American eel'
$ws.Cells.Item(353, 4).Value = 'Teleost Fish'
$ws.Cells.Item(354, 1).Value = '0e3aec812235602fac414c57ef969f1e'
$ws.Cells.Item(354, 2).Value = 'Peromyscus leucopus'
$ws.Cells.Item(354, 3).Value = 'Deer mouse'
$ws.Cells.Item(354, 4).Value = 'Mammal'
$ws.Cells.Item(355, 1).Value = '35a8484c22fbf1df676003af6ec52a29'
$ws.Cells.Item(355, 2).Value = 'Pseudopleuronectes americanus or Myzopsetta ferruginea'
$ws.Cells.Item(355, 3).Value = 'Winter or Yellowtail flounder'
$ws.Cells.Item(355, 4).Value = 'Teleost Fish'
$ws.Cells.Item(356, 1).Value = '841570d820eaab46bb7c7b3b7db8fba3'
$ws.Cells.Item(356, 2).Value = 'Ctenogobius boleosoma'
$ws.Cells.Item(356, 3).Value = 'Darter goby'
$ws.Cells.Item(356, 4).Value = 'Teleost Fish'
$ws.Cells.Item(357, 1).Value = 'b45f5a79f080475521114b63aa9bd7d1'
$ws.Cells.Item(357, 2).Value = 'Thunnus sp'
$ws.Cells.Item(357, 3).Value = 'Tuna sp'
$ws.Cells.Item(357, 4).Value = 'Teleost Fish'
$ws.Cells.Item(358, 1).Value = '5d6bb2fed75f92f01e645ffc80d17d36'
$ws.Cells.Item(358, 2).Value = 'Anchoa mitchilli'
$ws.Cells.Item(358, 3).Value = 'Bay anchovy'
$ws.Cells.Item(358, 4).Value = 'Teleost Fish'
$ws.Cells.Item(359, 1).Value = 'cc23248fc9f0058810041c6090c99461'
$ws.Cells.Item(359, 2).Value = 'Unassigned'
$ws.Cells.Item(359, 3).Value = 'Unassigned'
$ws.Cells.Item(359, 4).Value = 'Unassigned'
$ws.Cells.Item(364, 1).Value = '02af37069efe14191838e3ded56eba42'
$ws.Cells.Item(364, 2).Value = 'Unassigned'
$ws.Cells.Item(364, 3).Value = 'Unassigned'
$ws.Cells.Item(364, 4).Value = 'Unassigned'
$ws.Cells.Item(365, 1).Value = '0dc4976a75f5949215d7cf1f1a2994f6'
$ws.Cells.Item(365, 2).Value = 'Unassigned'
$ws.Cells.Item(365, 3).Value = 'Unassigned'
$ws.Cells.Item(365, 4).Value = 'Unassigned'
$ws.Cells.Item(366, 1).Value = 'dedc23ac2bedc9a49f0d9d5f2092f08a'
$ws.Cells.Item(366, 2).Value = 'Sternotherus carinatus'
$ws.Cells.Item(366, 3).Value = 'Razor-backed musk turtle'
$ws.Cells.Item(366, 4).Value = 'Reptile'
$ws.Cells.Item(368, 1).Value = '072eff78dfd4d6ee7d89b57c63c82827'
$ws.Cells.Item(368, 2).Value = 'Cottidae sp'
$ws.Cells.Item(368, 3).Value = 'Sculpins'
$ws.Cells.Item(368, 4).Value = 'Teleost Fish'
$ws.Cells.Item(369, 1).Value = 'f3ca4711d94796503be4e5bbb8a6705b'
$ws.Cells.Item(369, 2).Value = 'Homo sapiens'
$ws.Cells.Item(369, 3).Value = 'Human'
$ws.Cells.Item(369, 4).Value = 'Human'
$ws.Cells.Item(376, 1).Value = '66ce172600ad7ff5f35d5dfc0bab87d3'
$ws.Cells.Item(376, 2).Value = 'Homo sapiens'
$ws.Cells.Item(376, 3).Value = 'Human'
$ws.Cells.Item(376, 4).Value = 'Human'
$ws.Cells.Item(377, 1).Value = '9b1cbc97eff7dc6af2c100d85526140f'
$ws.Cells.Item(377, 2).Value = 'Homo sapiens'
$ws.Cells.Item(377, 3).Value = 'Human'
$ws.Cells.Item(377, 4).Value = 'Human'
$ws.Cells.Item(380, 1).Value = '798cab8e6a1a556a317f78cfa6bab8ac'
$ws.Cells.Item(380, 2).Value = 'Felis catus'
$ws.Cells.Item(380, 3).Value = 'Cat'
$ws.Cells.Item(380, 4).Value = 'Livestock'
$ws.Cells.Item(381, 1).Value = '97444d2388851f1d71afeb95125b4898'
$ws.Cells.Item(381, 2).Value = 'Unassigned'
$ws.Cells.Item(381, 3).Value = 'Unassigned'
$ws.Cells.Item(381, 4).Value = 'Unassigned'
$ws.Cells.Item(387, 1).Value = '41714252fc55bbced79f657eb2b8805a'
$ws.Cells.Item(387, 2).Value = 'Passer domesticus'
$ws.Cells.Item(387, 3).Value = 'House sparrow'
$ws.Cells.Item(387, 4).Value = 'Bird'
$ws.Cells.Item(389, 1).Value = '3b578403acdc73dd077d282c96f9541f'
$ws.Cells.Item(389, 2).Value = 'Ophidion marginatum'
$ws.Cells.Item(389, 3).Value = 'Striped cusk-eel'
$ws.Cells.Item(389, 4).Value = 'Teleost Fish'
$ws.Cells.Item(393, 1).Value = 'd1af29b8548fccd9f3eada2b18f0eac9'
$ws.Cells.Item(393, 2).Value = 'Unassigned'
$ws.Cells.Item(393, 3).Value = 'Unassigned'
$ws.Cells.Item(393, 4).Value = 'Unassigned'
$ws.Cells.Item(394, 1).Value = '803a43fa7cb74bb51f36ab2949523bf2'
$ws.Cells.Item(394, 2).Value = 'Homo sapiens'
$ws.Cells.Item(394, 3).Value = 'Human'
$ws.Cells.Item(394, 4).Value = 'Human'
$ws.Cells.Item(397, 1).Value = '812ed0386e2a4869a21da5634665548d'
$ws.Cells.Item(397, 2).Value = 'Caranx hippos'
$ws.Cells.Item(397, 3).Value = 'Crevalle jack'
$ws.Cells.Item(397, 4).Value = 'Teleost Fish'
$ws.Cells.Item(404, 1).Value = 'a6bf1361741b5eda21b4d05f18f04a90'
$ws.Cells.Item(404, 2).Value = 'Unassigned'
$ws.Cells.Item(404, 3).Value = 'Unassigned'
$ws.Cells.Item(404, 4).Value = 'Unassigned'
$ws.Cells.Item(405, 1).Value = '2d0506c060ee125f6608b52f22e598b1'
$ws.Cells.Item(405, 2).Value = 'Homo sapiens'
$ws.Cells.Item(405, 3).Value = 'Human'
$ws.Cells.Item(405, 4).Value = 'Human'
$ws.Cells.Item(407, 1).Value = '9e218ddde05826daea9943de26124674'
$ws.Cells.Item(407, 2).Value = 'Sturnus vulgaris'
$ws.Cells.Item(407, 3).Value = 'Common starling'
$ws.Cells.Item(407, 4).Value = 'Bird'
$ws.Cells.Item(408, 1).Value = 'f6b3a673a06591a2d5a8936584e64754'
$ws.Cells.Item(408, 2).Value = 'Tautogolabrus adspersus'
$ws.Cells.Item(408, 3).Value = 'Cunner'
$ws.Cells.Item(408, 4).Value = 'Teleost Fish'
$ws.Cells.Item(409, 1).Value = '75fedd0fa34e3ac2514601e68b613736'
$ws.Cells.Item(409, 2).Value = 'Strongylura marina'
$ws.Cells.Item(409, 3).Value = 'Atlantic needlefish'
$ws.Cells.Item(409, 4).Value = 'Teleost Fish'
$ws.Cells.Item(411, 1).Value = '037bd6992d173dfbcd22d76af622fa5b'
$ws.Cells.Item(411, 2).Value = 'Malaclemys terrapin'
$ws.Cells.Item(411, 3).Value = 'Diamondback terrapin'
$ws.Cells.Item(411, 4).Value = 'Reptile'
$ws.Cells.Item(412, 1).Value = '08bd987bd944513cc896ab3b3c3eed38'
$ws.Cells.Item(412, 2).Value = 'Hippoglossina oblonga'
$ws.Cells.Item(412, 3).Value = 'Fourspot flounder'
$ws.Cells.Item(412, 4).Value = 'Teleost Fish'
$ws.Cells.Item(413, 1).Value = '0437f2363acf453d6291ceda4abba683'
$ws.Cells.Item(413, 2).Value = 'Ameiurus nebulosus'
$ws.Cells.Item(413, 3).Value = 'Brown bullhead'
$ws.Cells.Item(413, 4).Value = 'Teleost Fish'
$ws.Cells.Item(414, 1).Value = '956d5064b9d6c222e19d75e231925e18'
$ws.Cells.Item(414, 2).Value = 'Homo sapiens'
$ws.Cells.Item(414, 3).Value = 'Human'
$ws.Cells.Item(414, 4).Value = 'Human'
$ws.Cells.Item(415, 1).Value = 'd964add43fe0c3212cbe19a066dc2a13'
$ws.Cells.Item(415, 2).Value = 'Unassigned'
$ws.Cells.Item(415, 3).Value = 'Unassigned'
$ws.Cells.Item(415, 4).Value = 'Unassigned'
$ws.Cells.Item(416, 1).Value = '7f400300a06f165c23af04aa4e4c790c'
$ws.Cells.Item(416, 2).Value = 'Enchelyopus cimbrius'
$ws.Cells.Item(416, 3).Value = 'Fourbeard rockling'
$ws.Cells.Item(416, 4).Value = 'Teleost Fish'
$ws.Cells.Item(417, 1).Value = 'ced5f183dc83bf9a3831984cefa3b3a1'
$ws.Cells.Item(417, 2).Value = 'Homo sapiens'
$ws.Cells.Item(417, 3).Value = 'Human'
$ws.Cells.Item(417, 4).Value = 'Human'
$ws.Cells.Item(418, 1).Value = '29ae99676d91135240ab43f0184c5909'
$ws.Cells.Item(418, 2).Value = 'Unassigned'
$ws.Cells.Item(418, 3).Value = 'Unassigned'
$ws.Cells.Item(418, 4).Value = 'Unassigned'
$ws.Cells.Item(419, 1).Value = 'df6846357baa6a9fecd66b4a1ba513a8'
$ws.Cells.Item(419, 2).Value = 'Canis lupus'
$ws.Cells.Item(419, 3).Value = 'Dog'
$ws.Cells.Item(419, 4).Value = 'Livestock'
$ws.Cells.Item(420, 1).Value = '39c09623e77e6cb1f69a264089e6256c'
$ws.Cells.Item(420, 2).Value = 'Blarina brevicauda'
$ws.Cells.Item(420, 3).Value = 'Northern short tailed shrew'
$ws.Cells.Item(420, 4).Value = 'Mammal'
